$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the fun and exciting features of Book of Helios slot game. Play for free and win big with the Book of the Sun symbol and free spins round.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Near the end of the document: remove the duplicate bold
#    "Play Book of Helios Slot Game for Free - Review" paragraph, and
#    replace the text of the following italic paragraph with the DALL-E
#    image-prompt text (keeping its italic formatting).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Text -eq "Play Book of Helios Slot Game for Free - Review`r") -and
        ($p.Style.NameLocal -ne "Heading 1")) {
        $p.Range.Delete() | Out-Null
        break
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastRange.Text = 'DALLE, please create a cartoon-style feature image for the game "Book of Helios" that features a happy Maya warrior with glasses. The image should be captivating and exciting, with the Maya warrior shown holding a copy of the book with Helios on the cover. The warrior should have a big smile and be surrounded by golden rays of sunlight. The image should also include other symbols from the game, such as the Book of Helios symbol and the expandable symbol. Make sure the overall color scheme is bright and eye-catching, and that the image is of high-quality. Thanks!'
